$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.448.93'
$ws.Range('E2').Value = '  -0.72%  '

$ws.Range('D3').Value = '2.518.65'
$ws.Range('E3').Value = '  -1.10%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '311.82'
$ws.Range('E5').Value = '  +1.17%  '

$ws.Range('D6').Value = '98.79'
$ws.Range('E6').Value = '  -3.55%  '

$ws.Range('E7').Value = '  -1.45%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  -3.17%  '

$ws.Range('D10').Value = '35.16'
$ws.Range('E10').Value = '  -3.48%  '

$ws.Range('D11').Value = '0.0801'
$ws.Range('E11').Value = '  -1.06%  '

$ws.Range('D12').Value = '0.109'
$ws.Range('E12').Value = '  +0.38%  '

$ws.Range('D13').Value = '7.20'
$ws.Range('E13').Value = '  -2.84%  '

$ws.Range('D14').Value = '2.909.91'
$ws.Range('E14').Value = '  -1.00%  '

$ws.Range('D15').Value = '15.27'
$ws.Range('E15').Value = '  -3.85%  '

$ws.Range('D16').Value = '2.516.79'
$ws.Range('E16').Value = '  -0.86%  '

$ws.Range('E17').Value = '  -3.89%  '

$ws.Range('D18').Value = '42.461.92'
$ws.Range('E18').Value = '  -0.79%  '

$ws.Range('D19').Value = '6.59'
$ws.Range('E19').Value = '  -2.99%  '

$ws.Range('D20').Value = '0.0₃0944'
$ws.Range('E20').Value = '  -1.36%  '

$ws.Range('D21').Value = '12.05'
$ws.Range('E21').Value = '  -2.85%  '

$ws.Range('D22').Value = '69.29'
$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('D23').Value = '239.90'
$ws.Range('E23').Value = '  -2.46%  '

$ws.Range('E24').Value = '  -2.02%  '

$ws.Range('E25').Value = '  -3.35%  '

$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('D27').Value = '25.26'
$ws.Range('E27').Value = '  -4.96%  '

$ws.Range('D28').Value = '2.26'
$ws.Range('E28').Value = '  -2.81%  '

$ws.Range('E29').Value = '  -1.18%  '

$ws.Range('D30').Value = '38.25'
$ws.Range('E30').Value = '  -6.65%  '

$ws.Range('D31').Value = '5.83'
$ws.Range('E31').Value = '  +1.65%  '

$ws.Range('D32').Value = '156.71'
$ws.Range('E32').Value = '  +0.29%  '

$ws.Range('D33').Value = '2.72'
$ws.Range('E33').Value = '  +3.85%  '

$ws.Range('D34').Value = '2.68'
$ws.Range('E34').Value = '  +1.58%  '

$ws.Range('D35').Value = '0.0785'
$ws.Range('E35').Value = '  -2.28%  '

$ws.Range('D36').Value = '3.13'
$ws.Range('E36').Value = '  -4.61%  '

$ws.Range('D37').Value = '1.95'
$ws.Range('E37').Value = '  -6.65%  '

$ws.Range('D38').Value = '17.45'
$ws.Range('E38').Value = '  -4.46%  '

$ws.Range('D39').Value = '0.108'
$ws.Range('E39').Value = '  -3.03%  '

$ws.Range('E40').Value = '  -1.05%  '

$ws.Range('D41').Value = '4.12'
$ws.Range('E41').Value = '  -2.21%  '

$ws.Range('D42').Value = '21.75'
$ws.Range('E42').Value = '  -2.91%  '

$ws.Range('E43').Value = '  +0.33%  '

$ws.Range('D44').Value = '3.26'
$ws.Range('E44').Value = '  -1.03%  '

$ws.Range('D45').Value = '0.0296'
$ws.Range('E45').Value = '  -1.38%  '

$ws.Range('D46').Value = '1.995.43'
$ws.Range('E46').Value = '  +0.63%  '

$ws.Range('D47').Value = '9.09'
$ws.Range('E47').Value = '  +1.06%  '

$ws.Range('D48').Value = '2.761.20'
$ws.Range('E48').Value = '  -1.17%  '

$ws.Range('E49').Value = '  -2.49%  '

$ws.Range('D50').Value = '78.59'
$ws.Range('E50').Value = '  -3.47%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '100.46'
$ws.Range('E51').Value = '  -1.48%  '
